# Localize the SharePoint "document" content-type schema (customXml/item1.xml)
# from Japanese to English display strings, bump its versionID/fieldsID, and
# re-mint the GUID of the companion datastore item (customXml/itemProps1.xml).
#
# This mirrors what "Merge Conflict resolution" / "Attempt to fix merge
# issues" commits typically do to the hidden SharePoint CustomXMLParts that
# PowerPoint carries around when a deck was checked in/out of a document
# library: the content-type schema part (identified by its root element
# <ct:contentTypeSchema .../>) gets its localized ma:* attributes swapped for
# the English originals, and the sibling itemProps part gets a fresh
# ds:itemID GUID.

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

# --- 1) customXml/item1.xml: the ct:contentTypeSchema part ------------------
# Locate it defensively by namespace/content rather than assuming a fixed
# index, since CustomXMLParts ordering is not guaranteed.
$schemaPart = $null
try {
    $matches = $parts.SelectByNamespace("http://schemas.microsoft.com/office/2006/metadata/contentType")
    if ($null -ne $matches -and $matches.Count -ge 1) {
        $schemaPart = $matches.Item(1)
    }
} catch {
    $schemaPart = $null
}

if ($null -eq $schemaPart) {
    for ($i = 1; $i -le $parts.Count; $i++) {
        $candidate = $parts.Item($i)
        if ($candidate.XML -like "*ct:contentTypeSchema*") {
            $schemaPart = $candidate
            break
        }
    }
}

if ($null -ne $schemaPart) {
    $xml = $schemaPart.XML

    # Localized content-type metadata -> English.
    $xml = $xml.Replace('ma:contentTypeName="ドキュメント"', 'ma:contentTypeName="Document"')
    $xml = $xml.Replace('ma:contentTypeDescription="新しいドキュメントを作成します。"', 'ma:contentTypeDescription="Create a new document."')

    # Re-minted version/fields stamps.
    $xml = $xml.Replace('ma:versionID="a9cab35011a557c1232e9e1918db7064"', 'ma:versionID="d0e002fabf17cb2440d8e9a473d3a41c"')
    $xml = $xml.Replace('ma:fieldsID="36c473bbc383ceb924bb8d2cdd9a2de6"', 'ma:fieldsID="e4cec627508c1f1ba247db94416ea198"')

    # Localized field display names -> English.
    $xml = $xml.Replace('ma:displayName="画像タグ"', 'ma:displayName="Image Tags"')
    $xml = $xml.Replace('ma:displayName="コンテンツ タイプ"', 'ma:displayName="Content Type"')
    $xml = $xml.Replace('ma:displayName="タイトル"', 'ma:displayName="Title"')

    $schemaPart.XML = $xml
}

# --- 2) customXml/itemProps1.xml: the companion datastore item GUID --------
$propsPart = $null
try {
    $matches2 = $parts.SelectByNamespace("http://schemas.openxmlformats.org/officeDocument/2006/customXml")
    if ($null -ne $matches2 -and $matches2.Count -ge 1) {
        $propsPart = $matches2.Item(1)
    }
} catch {
    $propsPart = $null
}

if ($null -eq $propsPart -and $null -ne $schemaPart) {
    # itemProps1.xml ships alongside item1.xml; re-use the same collection
    # lookup by scanning for the datastoreItem root if namespace lookup
    # above didn't resolve it.
    for ($i = 1; $i -le $parts.Count; $i++) {
        $candidate = $parts.Item($i)
        if ($candidate.XML -like "*ds:datastoreItem*8C2F49B3-9B4A-4186-81FD-94B9DADE14AD*") {
            $propsPart = $candidate
            break
        }
    }
}

if ($null -ne $propsPart) {
    $propsXml = $propsPart.XML
    $propsXml = $propsXml.Replace('{8C2F49B3-9B4A-4186-81FD-94B9DADE14AD}', '{F6D8E696-A6C0-4966-8CDE-8D4A3B2EC329}')
    $propsPart.XML = $propsXml
}

Write-Output "Updated customXml content-type schema and datastore item GUID."
